# "Change font in dfm_weight()" -- quanteda cheatsheet
#
# On slide 1, inside the "Weight or smooth the feature frequencies" code
# sample, the line:
#     dfm_weight(x, scheme = "prop") | dfm_smooth(x, smoothing = 0.5)
# has its "(x, scheme = "prop") | " segment re-run/re-fonted to explicitly
# carry the Monaco typeface (it previously relied on inherited formatting
# for part of that span). This locates the shape/run by text rather than a
# hard-coded index so it keeps working if shape ordering ever shifts.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Find the shape that contains the dfm_weight() code sample.
$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.TextRange.Text -like "*dfm_weight(x, scheme*") {
            $target = $shp
            break
        }
    }
}

$tr = $target.TextFrame.TextRange
$full = $tr.Text

# Locate "(x, scheme = "prop") | " right after "dfm_weight" (0-based offset).
$anchor = $full.IndexOf("dfm_weight(x, scheme")
$start0 = $anchor + "dfm_weight".Length

# ---- Split "(x, " into "(x" + ", " ---------------------------------------
# The comma+space becomes its own run carrying an explicit Monaco font.
$commaRun = $tr.Characters($start0 + 3, 2)
$commaRun.Font.Name = "Monaco"

# ---- Give the "scheme" run an explicit Monaco font -----------------------
$schemeRun = $tr.Characters($start0 + 5, 6)
$schemeRun.Font.Name = "Monaco"

# ---- Re-apply Monaco to ' = "prop")<nbsp>' --------------------------------
$propRun = $tr.Characters($start0 + 11, 11)
$propRun.Font.Name = "Monaco"

# ---- Split "| " off into its own run (keeps its inherited Monaco look) ---
$pipeRun = $tr.Characters($start0 + 22, 2)
$pipeRun.Font.Size = $pipeRun.Font.Size
